$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 288, shifting existing rows 288:369 down to 289:370
$ws.Rows("288:288").Insert()

# Populate the newly inserted row 288 with the new weekly record
$ws.Cells.Item(288, 1).Value = 4
$ws.Cells.Item(288, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(288, 3).Value = "Los Lagos"
$ws.Cells.Item(288, 4).Value = 44841
$ws.Cells.Item(288, 5).Value = 10
$ws.Cells.Item(288, 6).Value = 100112045
$ws.Cells.Item(288, 7).Value = "Zapallo"
$ws.Cells.Item(288, 8).Value = "Paine"
$ws.Cells.Item(288, 9).Value = "1a (guarda)"
$ws.Cells.Item(288, 10).Value = 1200
$ws.Cells.Item(288, 11).Value = 750
$ws.Cells.Item(288, 12).Value = 750
$ws.Cells.Item(288, 13).Value = 750
$ws.Cells.Item(288, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(288, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(288, 16).Value = 750
$ws.Cells.Item(288, 17).Value = 1
$ws.Cells.Item(288, 18).Value = "Hortaliza"
